# Add a new "ProjectGroup" column to the worksheet / Table1, populate it
# with sequential values, remove the stray ProjectDependency value that
# used to live in C3, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table1 currently spans A1:C5 (ProjectID, ProjectName, ProjectDependency).
$tbl = $ws.ListObjects.Item(1)

# The lone "ProjectDependency" data point (C3) is being dropped.
$ws.Range("C3").ClearContents() | Out-Null

# Grow the table by one column; this shifts ref/autoFilter to A1:D5 and
# adds a 4th <tableColumn> automatically.
$newColumn = $tbl.ListColumns.Add()

# Naming the column through the header cell keeps the shared-string table
# and the table XML's column name in sync.
$ws.Range("D1").Value = "ProjectGroup"
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 4

# Match the saved selection state from the edited file.
$ws.Range("E7").Select() | Out-Null
